$d = $word.ActiveDocument

# Replace each unique multiplication-fact cell value with its new value.
# Parameters: FindText, MatchCase, MatchWholeWord, MatchWildcards,
#             MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#             Format, ReplaceWith, Replace (2 = wdReplaceAll)

$d.Content.Find.Execute("528×4=2112", $true, $false, $false, $false, $false, $true, 1, $false, "480×7=3360", 2) | Out-Null
$d.Content.Find.Execute("687×9=6183", $true, $false, $false, $false, $false, $true, 1, $false, "221×8=1768", 2) | Out-Null
$d.Content.Find.Execute("340×3=1020", $true, $false, $false, $false, $false, $true, 1, $false, "590×5=2950", 2) | Out-Null
$d.Content.Find.Execute("755×4=3020", $true, $false, $false, $false, $false, $true, 1, $false, "706×3=2118", 2) | Out-Null
$d.Content.Find.Execute("132×8=1056", $true, $false, $false, $false, $false, $true, 1, $false, "459×8=3672", 2) | Out-Null
$d.Content.Find.Execute("457×6=2742", $true, $false, $false, $false, $false, $true, 1, $false, "395×6=2370", 2) | Out-Null
$d.Content.Find.Execute("753×8=6024", $true, $false, $false, $false, $false, $true, 1, $false, "726×2=1452", 2) | Out-Null
$d.Content.Find.Execute("687×7=4809", $true, $false, $false, $false, $false, $true, 1, $false, "924×9=8316", 2) | Out-Null
$d.Content.Find.Execute("674×4=2696", $true, $false, $false, $false, $false, $true, 1, $false, "340×3=1020", 2) | Out-Null
$d.Content.Find.Execute("134×5=670", $true, $false, $false, $false, $false, $true, 1, $false, "817×7=5719", 2) | Out-Null
$d.Content.Find.Execute("340×8=2720", $true, $false, $false, $false, $false, $true, 1, $false, "332×5=1660", 2) | Out-Null
$d.Content.Find.Execute("764×7=5348", $true, $false, $false, $false, $false, $true, 1, $false, "869×9=7821", 2) | Out-Null
$d.Content.Find.Execute("939×7=6573", $true, $false, $false, $false, $false, $true, 1, $false, "399×4=1596", 2) | Out-Null
$d.Content.Find.Execute("126×7=882", $true, $false, $false, $false, $false, $true, 1, $false, "836×4=3344", 2) | Out-Null
$d.Content.Find.Execute("556×9=5004", $true, $false, $false, $false, $false, $true, 1, $false, "612×6=3672", 2) | Out-Null
$d.Content.Find.Execute("771×2=1542", $true, $false, $false, $false, $false, $true, 1, $false, "629×9=5661", 2) | Out-Null
$d.Content.Find.Execute("524×4=2096", $true, $false, $false, $false, $false, $true, 1, $false, "152×4=608", 2) | Out-Null
$d.Content.Find.Execute("373×8=2984", $true, $false, $false, $false, $false, $true, 1, $false, "951×5=4755", 2) | Out-Null
$d.Content.Find.Execute("981×9=8829", $true, $false, $false, $false, $false, $true, 1, $false, "879×7=6153", 2) | Out-Null
$d.Content.Find.Execute("252×8=2016", $true, $false, $false, $false, $false, $true, 1, $false, "646×5=3230", 2) | Out-Null
$d.Content.Find.Execute("537×6=3222", $true, $false, $false, $false, $false, $true, 1, $false, "338×4=1352", 2) | Out-Null
$d.Content.Find.Execute("399×9=3591", $true, $false, $false, $false, $false, $true, 1, $false, "202×7=1414", 2) | Out-Null
$d.Content.Find.Execute("221×9=1989", $true, $false, $false, $false, $false, $true, 1, $false, "321×6=1926", 2) | Out-Null
$d.Content.Find.Execute("272×9=2448", $true, $false, $false, $false, $false, $true, 1, $false, "562×5=2810", 2) | Out-Null
$d.Content.Find.Execute("161×3=483", $true, $false, $false, $false, $false, $true, 1, $false, "847×2=1694", 2) | Out-Null
